# Fruta / hortaliza, semanal
# Inserts the latest weekly report (2 new rows: "Pintón" and "Primera Pintón")
# at the top of the Platano data block (row 382), pushing the rest of the
# historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 382; this shifts every existing row
# (382..497) down to (384..499) and automatically grows the used range to
# A1:T499, matching the new <dimension> in the target workbook.
$ws.Rows("382:383").Insert()

# New row 382: "Pintón" quality entry for the new reporting date.
$ws.Cells.Item(382, 1).Value  = 7
$ws.Cells.Item(382, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(382, 3).Value  = "Ñuble"
$ws.Cells.Item(382, 4).Value  = 44588
$ws.Cells.Item(382, 5).Value  = 16
$ws.Cells.Item(382, 6).Value  = "Fruta"
$ws.Cells.Item(382, 7).Value  = 100108
$ws.Cells.Item(382, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(382, 9).Value  = 100108006
$ws.Cells.Item(382, 10).Value = "Plátano"
$ws.Cells.Item(382, 11).Value = "Sin especificar"
$ws.Cells.Item(382, 12).Value = "Pintón"
$ws.Cells.Item(382, 13).Value = 300
$ws.Cells.Item(382, 14).Value = 11000
$ws.Cells.Item(382, 15).Value = 11000
$ws.Cells.Item(382, 16).Value = 11000
$ws.Cells.Item(382, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(382, 18).Value = "Ecuador"
$ws.Cells.Item(382, 19).Value = 550
$ws.Cells.Item(382, 20).Value = 20

# New row 383: "Primera Pintón" quality entry for the same reporting date.
$ws.Cells.Item(383, 1).Value  = 7
$ws.Cells.Item(383, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(383, 3).Value  = "Ñuble"
$ws.Cells.Item(383, 4).Value  = 44588
$ws.Cells.Item(383, 5).Value  = 16
$ws.Cells.Item(383, 6).Value  = "Fruta"
$ws.Cells.Item(383, 7).Value  = 100108
$ws.Cells.Item(383, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(383, 9).Value  = 100108006
$ws.Cells.Item(383, 10).Value = "Plátano"
$ws.Cells.Item(383, 11).Value = "Sin especificar"
$ws.Cells.Item(383, 12).Value = "Primera Pintón"
$ws.Cells.Item(383, 13).Value = 600
$ws.Cells.Item(383, 14).Value = 12000
$ws.Cells.Item(383, 15).Value = 13000
$ws.Cells.Item(383, 16).Value = 12500
$ws.Cells.Item(383, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(383, 18).Value = "Ecuador"
$ws.Cells.Item(383, 19).Value = 625
$ws.Cells.Item(383, 20).Value = 20
